# Hortaliza, Terminal Hortofrutícola Agro Chillán - Perejil
# Weekly price update: a new week of data (date 45239) is added for the
# "Primera" and "Segunda" quality rows of Ñuble, which is implemented by
# inserting two fresh rows right above the former last two rows (copies
# of the data that is about to be overwritten), then updating the
# original rows 143/144 in place with the new week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert two blank rows before row 145 -- this pushes the former
#    row 145 down to row 147 and leaves two empty rows at 145 and 146.
$ws.Range("A145:A146").EntireRow.Insert()

# 2) Preserve the previous content of (what were) rows 143 and 144 by
#    copying them down into the freshly inserted rows 145 and 146,
#    before they get overwritten with the new week's values below.
$ws.Range("A143:R143").Copy()
$ws.Range("A145").PasteSpecial()

$ws.Range("A144:R144").Copy()
$ws.Range("A146").PasteSpecial()

$excel.CutCopyMode = $false

# 3) Update row 143 (quality "Primera") with the new week's data.
$ws.Cells.Item(143, 4).Value = 45239
$ws.Cells.Item(143, 10).Value = 300
$ws.Cells.Item(143, 11).Value = 2000
$ws.Cells.Item(143, 12).Value = 2000
$ws.Cells.Item(143, 13).Value = 2000
$ws.Cells.Item(143, 16).Value = 2000

# 4) Update row 144 (now quality "Segunda", region "Región de Ñuble")
#    with the new week's data.
$ws.Cells.Item(144, 4).Value = 45239
$ws.Cells.Item(144, 9).Value = "Segunda"
$ws.Cells.Item(144, 10).Value = 200
$ws.Cells.Item(144, 11).Value = 1500
$ws.Cells.Item(144, 12).Value = 1500
$ws.Cells.Item(144, 13).Value = 1500
$ws.Cells.Item(144, 15).Value = "Región de Ñuble"
$ws.Cells.Item(144, 16).Value = 1500
